$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new attendance column (AV) for the next training day (2025-09-17,
# serial 45917), following the existing AU column (2025-09-16, serial 45916).
# Values mirror each player's previous-day (AU) status, except row 21
# (Amir Kherrab) who is now "P" (Présent) instead of "B" (Blessure).

# 1) Write the new date header.
$ws.Range("AV1").Value = 45917

# 2) Write each player's attendance mark for the new day.
$ws.Range("AV2").Value  = "P"
$ws.Range("AV3").Value  = "P"
$ws.Range("AV4").Value  = "P"
$ws.Range("AV5").Value  = "B"
$ws.Range("AV6").Value  = "B"
$ws.Range("AV7").Value  = "P"
$ws.Range("AV8").Value  = "P"
$ws.Range("AV9").Value  = "P"
$ws.Range("AV10").Value = "B"
$ws.Range("AV11").Value = "P"
$ws.Range("AV12").Value = "P"
$ws.Range("AV13").Value = "B"
$ws.Range("AV14").Value = "P"
$ws.Range("AV15").Value = "P"
$ws.Range("AV16").Value = "P"
$ws.Range("AV17").Value = "P"
$ws.Range("AV18").Value = "P"
$ws.Range("AV19").Value = "P"
$ws.Range("AV20").Value = "P"
$ws.Range("AV21").Value = "P"
$ws.Range("AV22").Value = "P"
$ws.Range("AV23").Value = "P"
$ws.Range("AV24").Value = "P"
$ws.Range("AV25").Value = "P"
$ws.Range("AV26").Value = "P"
$ws.Range("AV27").Value = "P"
$ws.Range("AV28").Value = "P"
$ws.Range("AV29").Value = "P"

# 3) Copy the formatting (number format / alignment / style) from the
#    previous date column (AU) onto the new column (AV) so the new cells
#    look like every other date column. Done after the values are set so
#    the dependent COUNTA/COUNTIF summary formulas pick up the new data.
$ws.Range("AU1:AU29").Copy()
$ws.Range("AV1:AV29").PasteSpecial(-4122)

# 4) Restore the previous selection behaviour (last cell the user clicked).
$ws.Range("AX24").Select()
